$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attack Effects")
$ws.Name = "Attacks"
